$wb = $excel.ActiveWorkbook

# --- Add the new "Goal 1" sheet before the existing "Goal 5" sheet ---
$ws0 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($ws0)
$newSheet.Name = "Goal 1"

# Re-fetch sheets by name to avoid stale COM anchors after Add()
$ws = $wb.Worksheets.Item("Goal 1")
$goal5 = $wb.Worksheets.Item("Goal 5")

# --- Cell values for "Goal 1" ---
$ws.Range("A1").Value = "World"
$ws.Range("B1").Value = "Min"
$ws.Range("C1").Value = "Avg"
$ws.Range("D1").Value = "Max"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 49260.081502045701
$ws.Range("C2").Value = 70447.732474781893
$ws.Range("D2").Value = 128596.92701965899
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 34887.1846051229
$ws.Range("C3").Value = 44184.556325680504
$ws.Range("D3").Value = 55952.754051457203
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 43984.346516044803
$ws.Range("C4").Value = 71461.934882409201
$ws.Range("D4").Value = 170129.46007863199
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 60097.9690427717
$ws.Range("C5").Value = 217504.94307903
$ws.Range("D5").Value = 692209.77471155103
$ws.Range("F5").Value = "Time"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 35056.049801858499
$ws.Range("C6").Value = 41933.495021554503
$ws.Range("D6").Value = 50785.410507758599
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 34952.697042174099
$ws.Range("C7").Value = 50276.5917783205
$ws.Range("D7").Value = 62050.285459474202
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 44008.751042916803
$ws.Range("C8").Value = 55034.7844003352
$ws.Range("D8").Value = 67762.5807125641
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 29908.4532416406
$ws.Range("C9").Value = 52759.893584238402
$ws.Range("D9").Value = 75567.226973988902
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 35091.061916969396
$ws.Range("C10").Value = 48009.464141524099
$ws.Range("D10").Value = 63424.394427009698
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 34918.595787494101
$ws.Range("C11").Value = 46168.599362492001
$ws.Range("D11").Value = 64545.752171613902
$ws.Range("A12").Value = "World"
$ws.Range("B12").Value = "Min"
$ws.Range("C12").Value = "Avg"
$ws.Range("D12").Value = "Max"
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 11
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 4.375
$ws.Range("D14").Value = 8
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 5.2857142857142803
$ws.Range("D15").Value = 7
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = 7.8181818181818103
$ws.Range("D16").Value = 9
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 3
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("F17").Value = "Sub-symbolic actions"
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 4.2857142857142803
$ws.Range("D18").Value = 5
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 4.5714285714285703
$ws.Range("D19").Value = 5
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 5
$ws.Range("C20").Value = 6.4166666666666599
$ws.Range("D20").Value = 8
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = 3.8571428571428501
$ws.Range("D21").Value = 5
$ws.Range("A22").Value = 9
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 3.71428571428571
$ws.Range("D22").Value = 5
$ws.Range("A23").Value = "World"
$ws.Range("C23").Value = "Avg"
$ws.Range("D23").Value = "Max"
$ws.Range("A24").Value = 0
$ws.Range("B24").Value = 21
$ws.Range("C24").Value = 31.636363636363601
$ws.Range("D24").Value = 42
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 15
$ws.Range("C25").Value = 19.375
$ws.Range("D25").Value = 28
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = 19
$ws.Range("C26").Value = 25.857142857142801
$ws.Range("D26").Value = 38
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = 32
$ws.Range("D27").Value = 40
$ws.Range("A28").Value = 4
$ws.Range("B28").Value = 15
$ws.Range("C28").Value = 18.285714285714199
$ws.Range("D28").Value = 22
$ws.Range("F28").Value = "Symbolic actions"
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = 15
$ws.Range("C29").Value = 22
$ws.Range("D29").Value = 28
$ws.Range("A30").Value = 6
$ws.Range("B30").Value = 19
$ws.Range("C30").Value = 22.857142857142801
$ws.Range("D30").Value = 26
$ws.Range("A31").Value = 7
$ws.Range("B31").Value = 16
$ws.Range("C31").Value = 26.1666666666666
$ws.Range("D31").Value = 34
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = 15
$ws.Range("C32").Value = 20.428571428571399
$ws.Range("D32").Value = 28
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = 15
$ws.Range("C33").Value = 19.571428571428498
$ws.Range("D33").Value = 28

# --- Conditional formatting (data bars) for "Goal 1" ---
$cf = $ws.Range("B2:B11").FormatConditions.AddDatabar()
$cf.Priority = 9
$cf = $ws.Range("B13:B21").FormatConditions.AddDatabar()
$cf.Priority = 8
$cf = $ws.Range("B24:B33").FormatConditions.AddDatabar()
$cf.Priority = 7
$cf = $ws.Range("C2:C11").FormatConditions.AddDatabar()
$cf.Priority = 6
$cf = $ws.Range("C13:C22").FormatConditions.AddDatabar()
$cf.Priority = 5
$cf = $ws.Range("C24:C33").FormatConditions.AddDatabar()
$cf.Priority = 4
$cf = $ws.Range("D2:D11").FormatConditions.AddDatabar()
$cf.Priority = 3
$cf = $ws.Range("D13:D22").FormatConditions.AddDatabar()
$cf.Priority = 2
$cf = $ws.Range("D24:D33").FormatConditions.AddDatabar()
$cf.Priority = 1

# --- Column widths for "Goal 1" (best achievable via ColumnWidth's character-unit rounding) ---
$ws.Range("B1:C1").ColumnWidth = 17.619791666666668
$ws.Range("D1").ColumnWidth = 15.346354166666666

# --- View/selection state ---
# Goal 5 loses the tab-selected flag and gets a new selection/range.
$goal5.Activate()
$goal5.Range("A1:I40").Select()

# Goal 1 becomes the active/selected tab with its own selection.
$ws.Activate()
$ws.Range("J42").Select()
